$wb = $excel.ActiveWorkbook

# Sheets that contain this data table: "展览" (Exhibition) and "全部类型" (All types)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1451
    $ws.Range("F7").Value = 8
    $ws.Range("F9").Value = 239
}
